$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold, centered, border) from an existing header cell
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-18
$values = @(
    @(7, 7),
    @(10, 10),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(4, 4),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(4, 4),
    @(3, 4),
    @(7, 7),
    @(6, 6),
    @(5, 5)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
